$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7720000147819519
$ws.Range("C2").Value = 0.8700000643730164
$ws.Range("D2").Value = 0.8210000395774841
$ws.Range("E2").Value = 0.7350000143051147

$ws.Range("B3").Value = 0.7710000276565552
$ws.Range("E3").Value = 0.7330000400543213

$ws.Range("B4").Value = 0.8347107172012329
$ws.Range("C4").Value = 0.8125
$ws.Range("D4").Value = 0.8909090757369995
$ws.Range("E4").Value = 0.7835820913314819

$ws.Range("B5").Value = 0.8264462351799011
$ws.Range("D5").Value = 0.9181817770004272
$ws.Range("E5").Value = 0.7985074520111084

$ws.Range("B6").Value = 0.7633674144744873
$ws.Range("C6").Value = 0.8709349036216736
$ws.Range("D6").Value = 0.8012820482254028
$ws.Range("E6").Value = 0.7274826765060425

$ws.Range("B7").Value = 0.7633674144744873
$ws.Range("C7").Value = 0.9999999403953552
$ws.Range("D7").Value = 0.8589743971824646
$ws.Range("E7").Value = 0.7228637337684631

$ws.Range("B8").Value = 0.8680000305175781
$ws.Range("D8").Value = 0.7730000615119934
$ws.Range("E8").Value = 0.8560000658035278

$ws.Range("B9").Value = 0.8790000677108765
$ws.Range("C9").Value = 0.984000027179718
$ws.Range("D9").Value = 0.7800000309944153
$ws.Range("E9").Value = 0.8660000562667847

$ws.Range("B10").Value = 1.090567111968994
$ws.Range("D10").Value = 1.073778629302979
$ws.Range("E10").Value = 1.145163059234619

$ws.Range("B11").Value = 1.089323759078979
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = 1.1527019739151
$ws.Range("E11").Value = 1.199325203895569

$ws.Range("B12").Value = 1.090567111968994
$ws.Range("D12").Value = 1.073778629302979
$ws.Range("E12").Value = 1.145163059234619
